$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-9
# from 45170 (2023-09-01) to 45174 (2023-09-05), keeping existing formatting.
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45174
}
